$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 21742858
$ws.Range("I70").Value = 2458.6365
$ws.Range("J70").Value = 41671556
$ws.Range("K70").Value = 7375.9095
$ws.Range("L70").Value = 125014668
$ws.Range("M70").Value = -7105.9095
$ws.Range("N70").Value = -125015208
$ws.Range("H73").Value = 21742858
$ws.Range("I73").Value = 2458.6365
$ws.Range("J73").Value = 41671556
$ws.Range("K73").Value = 7375.9095
$ws.Range("L73").Value = 125014668
$ws.Range("M73").Value = -6439.9095
$ws.Range("N73").Value = -125016540
$ws.Range("H92").Value = 47619824
$ws.Range("I92").Value = 50000668
$ws.Range("K92").Value = 50000668
$ws.Range("M92").Value = -49999420
$ws.Range("H98").Value = 1230.9062
$ws.Range("I98").Value = 1096.3
$ws.Range("K98").Value = 1096.3
$ws.Range("M98").Value = 401.7
$ws.Range("H112").Value = 1471.3529
$ws.Range("J112").Value = 1455.3541
$ws.Range("L112").Value = 4366.0623
$ws.Range("N112").Value = -6582.0623
$ws.Range("H122").Value = 1230.9062
$ws.Range("I122").Value = 1096.3
$ws.Range("K122").Value = 3288.9
$ws.Range("M122").Value = -838.8999999999996
$ws.Range("H132").Value = 113393.44
$ws.Range("I132").Value = 150304.16
$ws.Range("J132").Value = 17104.61
$ws.Range("K132").Value = 450912.48
$ws.Range("L132").Value = 51313.83
$ws.Range("M132").Value = -448382.48
$ws.Range("N132").Value = -56373.83
$ws.Range("H138").Value = 2668.4065
$ws.Range("J138").Value = 2942.0725
$ws.Range("L138").Value = 8826.217500000001
$ws.Range("N138").Value = -19106.2175

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13581.851
$ws.Range("I32").Value = 12808.1
$ws.Range("K32").Value = 12808.1
$ws.Range("M32").Value = -12521.1
$ws.Range("H61").Value = 8543.308000000001
$ws.Range("I61").Value = 8642.137000000001
$ws.Range("K61").Value = 8642.137000000001
$ws.Range("M61").Value = -8430.137000000001
$ws.Range("H63").Value = 4558.6
$ws.Range("I63").Value = 3400
$ws.Range("K63").Value = 3400
$ws.Range("M63").Value = -2714
$ws.Range("H66").Value = 4558.6
$ws.Range("I66").Value = 3400
$ws.Range("K66").Value = 17000
$ws.Range("M66").Value = -13568
$ws.Range("H74").Value = 27780052
$ws.Range("I74").Value = 83334240
$ws.Range("J74").Value = 2958
$ws.Range("K74").Value = 83334240
$ws.Range("L74").Value = 2958
$ws.Range("M74").Value = -83333366
$ws.Range("N74").Value = -4706
$ws.Range("H77").Value = 27780052
$ws.Range("I77").Value = 83334240
$ws.Range("J77").Value = 2958
$ws.Range("K77").Value = 416671200
$ws.Range("L77").Value = 14790
$ws.Range("M77").Value = -416666832
$ws.Range("N77").Value = -23526
$ws.Range("H110").Value = 1281352.1
$ws.Range("J110").Value = 18996.334
$ws.Range("L110").Value = 18996.334
$ws.Range("N110").Value = -23086.334
$ws.Range("H132").Value = 14134.4
$ws.Range("I132").Value = 18991.281
$ws.Range("K132").Value = 56973.84299999999
$ws.Range("M132").Value = -54443.84299999999
$ws.Range("H136").Value = 8543.308000000001
$ws.Range("I136").Value = 8642.137000000001
$ws.Range("K136").Value = 25926.411
$ws.Range("M136").Value = -23376.411

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4069.75
$ws.Range("I20").Value = 3751.7727
$ws.Range("J20").Value = 4769.3
$ws.Range("K20").Value = 3751.7727
$ws.Range("L20").Value = 4769.3
$ws.Range("M20").Value = -3504.7727
$ws.Range("N20").Value = -5263.3
$ws.Range("H86").Value = 1953.8462
$ws.Range("I86").Value = 1491.05
$ws.Range("K86").Value = 1491.05
$ws.Range("M86").Value = -368.05
$ws.Range("H89").Value = 1953.8462
$ws.Range("I89").Value = 1491.05
$ws.Range("K89").Value = 7455.25
$ws.Range("M89").Value = -1839.25
$ws.Range("H134").Value = 2557.25
$ws.Range("I134").Value = 2557.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7671.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5136.75
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 80262.664
$ws.Range("J141").Value = 94684.89
$ws.Range("L141").Value = 94684.89
$ws.Range("N141").Value = -105044.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 80.25
$ws.Range("I8").Value = 80.25
$ws.Range("K8").Value = 240.75
$ws.Range("M8").Value = -101.75
$ws.Range("H56").Value = 71435730
$ws.Range("I56").Value = 71435730
$ws.Range("K56").Value = 71435730
$ws.Range("M56").Value = -71435200
$ws.Range("H82").Value = 5479.125
$ws.Range("I82").Value = 3333
$ws.Range("J82").Value = 5785.7144
$ws.Range("K82").Value = 9999
$ws.Range("L82").Value = 17357.1432
$ws.Range("M82").Value = -9593
$ws.Range("N82").Value = -18169.1432
$ws.Range("H85").Value = 5479.125
$ws.Range("I85").Value = 3333
$ws.Range("J85").Value = 5785.7144
$ws.Range("K85").Value = 9999
$ws.Range("L85").Value = 17357.1432
$ws.Range("M85").Value = -8595
$ws.Range("N85").Value = -20165.1432
$ws.Range("H98").Value = 1414.4615
$ws.Range("J98").Value = 1451.0834
$ws.Range("L98").Value = 4353.2502
$ws.Range("N98").Value = -7349.2502
$ws.Range("H129").Value = 1119.875
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500
$ws.Range("H131").Value = 1896155.1
$ws.Range("I131").Value = 918.5
$ws.Range("J131").Value = 2068449.4
$ws.Range("K131").Value = 2755.5
$ws.Range("L131").Value = 6205348.199999999
$ws.Range("M131").Value = 2284.5
$ws.Range("N131").Value = -6215428.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1628493.6
$ws.Range("I70").Value = 2530123
$ws.Range("K70").Value = 2530123
$ws.Range("M70").Value = -2529853
$ws.Range("H73").Value = 1628493.6
$ws.Range("I73").Value = 2530123
$ws.Range("K73").Value = 2530123
$ws.Range("M73").Value = -2529187
$ws.Range("H80").Value = 1048148.44
$ws.Range("I80").Value = 1387189.6
$ws.Range("K80").Value = 1387189.6
$ws.Range("M80").Value = -1386191.6
$ws.Range("H83").Value = 1048148.44
$ws.Range("I83").Value = 1387189.6
$ws.Range("K83").Value = 6935948
$ws.Range("M83").Value = -6930956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 887.8182
$ws.Range("I22").Value = 817.2
$ws.Range("K22").Value = 817.2
$ws.Range("M22").Value = -522.2
$ws.Range("H27").Value = 887.8182
$ws.Range("I27").Value = 817.2
$ws.Range("K27").Value = 817.2
$ws.Range("M27").Value = -710.2
$ws.Range("H61").Value = 1567.5
$ws.Range("I61").Value = 1567.5
$ws.Range("K61").Value = 1567.5
$ws.Range("M61").Value = -1365.5
$ws.Range("H113").Value = 1567.5
$ws.Range("I113").Value = 1567.5
$ws.Range("K113").Value = 1567.5
$ws.Range("M113").Value = 602.5
$ws.Range("H122").Value = 3467.7646
$ws.Range("I122").Value = 2996.6785
$ws.Range("K122").Value = 8990.0355
$ws.Range("M122").Value = -6540.0355
$ws.Range("H136").Value = 4874.75
$ws.Range("I136").Value = 2166.1667
$ws.Range("K136").Value = 6498.500100000001
$ws.Range("M136").Value = -3948.500100000001
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 7162.2
$ws.Range("J29").Value = 8603.666999999999
$ws.Range("L29").Value = 8603.666999999999
$ws.Range("N29").Value = -9183.666999999999
$ws.Range("H62").Value = 5269814
$ws.Range("J62").Value = 12837.5
$ws.Range("L62").Value = 12837.5
$ws.Range("N62").Value = -14085.5
$ws.Range("H65").Value = 5269814
$ws.Range("J65").Value = 12837.5
$ws.Range("L65").Value = 64187.5
$ws.Range("N65").Value = -70427.5
$ws.Range("H81").Value = 2611652
$ws.Range("I81").Value = 3476370.2
$ws.Range("J81").Value = 17497.5
$ws.Range("K81").Value = 6952740.4
$ws.Range("L81").Value = 34995
$ws.Range("M81").Value = -6951679.4
$ws.Range("N81").Value = -37117
$ws.Range("H84").Value = 2611652
$ws.Range("I84").Value = 3476370.2
$ws.Range("J84").Value = 17497.5
$ws.Range("K84").Value = 34763702
$ws.Range("L84").Value = 174975
$ws.Range("M84").Value = -34758398
$ws.Range("N84").Value = -185583
$ws.Range("H136").Value = 5704.11
$ws.Range("I136").Value = 2320.0557
$ws.Range("K136").Value = 6960.1671
$ws.Range("M136").Value = -4410.1671
